$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 120
$ws.Range("I2").Value = 120
$ws.Range("K2").Value = 120
$ws.Range("M2").Value = -7
$ws.Range("H5").Value = 724.65
$ws.Range("I5").Value = 679.5333000000001
$ws.Range("J5").Value = 860
$ws.Range("K5").Value = 679.5333000000001
$ws.Range("L5").Value = 860
$ws.Range("M5").Value = -564.5333000000001
$ws.Range("N5").Value = -1090
$ws.Range("H9").Value = 411.88235
$ws.Range("I9").Value = 229.38461
$ws.Range("K9").Value = 229.38461
$ws.Range("M9").Value = -60.38461000000001
$ws.Range("H15").Value = 783129.5
$ws.Range("I15").Value = 783129.5
$ws.Range("K15").Value = 2349388.5
$ws.Range("M15").Value = -2349219.5
$ws.Range("H32").Value = 1964.6666
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H40").Value = 3963.5789
$ws.Range("J40").Value = 4023.8667
$ws.Range("L40").Value = 4023.8667
$ws.Range("N40").Value = -4373.8667
$ws.Range("H43").Value = 5923.2563
$ws.Range("I43").Value = 6141.0415
$ws.Range("J43").Value = 5574.8
$ws.Range("K43").Value = 6141.0415
$ws.Range("L43").Value = 5574.8
$ws.Range("M43").Value = -6072.0415
$ws.Range("N43").Value = -5712.8
$ws.Range("H51").Value = 7733.625
$ws.Range("J51").Value = 7838.4287
$ws.Range("L51").Value = 7838.4287
$ws.Range("N51").Value = -8806.4287
$ws.Range("H62").Value = 5351.7144
$ws.Range("I62").Value = 5502.6
$ws.Range("J62").Value = 4974.5
$ws.Range("K62").Value = 5502.6
$ws.Range("L62").Value = 4974.5
$ws.Range("M62").Value = -4878.6
$ws.Range("N62").Value = -6222.5
$ws.Range("H65").Value = 5351.7144
$ws.Range("I65").Value = 5502.6
$ws.Range("J65").Value = 4974.5
$ws.Range("K65").Value = 27513
$ws.Range("L65").Value = 24872.5
$ws.Range("M65").Value = -24393
$ws.Range("N65").Value = -31112.5
$ws.Range("H70").Value = 1251.76
$ws.Range("I70").Value = 1162.125
$ws.Range("K70").Value = 3486.375
$ws.Range("M70").Value = -3216.375
$ws.Range("H73").Value = 1251.76
$ws.Range("I73").Value = 1162.125
$ws.Range("K73").Value = 3486.375
$ws.Range("M73").Value = -2550.375
$ws.Range("H74").Value = 5823
$ws.Range("I74").Value = 5615.294
$ws.Range("K74").Value = 5615.294
$ws.Range("M74").Value = -4679.294
$ws.Range("H77").Value = 5823
$ws.Range("I77").Value = 5615.294
$ws.Range("K77").Value = 28076.47
$ws.Range("M77").Value = -23396.47
$ws.Range("H87").Value = 121416.664
$ws.Range("J87").Value = 129727.27
$ws.Range("L87").Value = 129727.27
$ws.Range("N87").Value = -132223.27
$ws.Range("H90").Value = 121416.664
$ws.Range("J90").Value = 129727.27
$ws.Range("L90").Value = 389181.81
$ws.Range("N90").Value = -401661.81
$ws.Range("H103").Value = 230
$ws.Range("J103").Value = 312.85715
$ws.Range("L103").Value = 938.5714499999999
$ws.Range("N103").Value = -2110.57145
$ws.Range("H107").Value = 957.0769
$ws.Range("I107").Value = 845.8
$ws.Range("J107").Value = 1328
$ws.Range("K107").Value = 845.8
$ws.Range("L107").Value = 1328
$ws.Range("M107").Value = 1074.2
$ws.Range("N107").Value = -5168
$ws.Range("H116").Value = 7765.6665
$ws.Range("I116").Value = 6453.273
$ws.Range("K116").Value = 6453.273
$ws.Range("M116").Value = -3011.273
$ws.Range("H125").Value = 84059.336
$ws.Range("J125").Value = 143702.72
$ws.Range("L125").Value = 1293324.48
$ws.Range("N125").Value = -1298244.48
$ws.Range("H132").Value = 3019.0344
$ws.Range("I132").Value = 3094.3333
$ws.Range("J132").Value = 2002.5
$ws.Range("K132").Value = 9282.999899999999
$ws.Range("L132").Value = 6007.5
$ws.Range("M132").Value = -6752.999899999999
$ws.Range("N132").Value = -11067.5
$ws.Range("H134").Value = 84999.75
$ws.Range("J134").Value = 84999.75
$ws.Range("L134").Value = 84999.75
$ws.Range("N134").Value = -95139.75
$ws.Range("H135").Value = 2173.7856
$ws.Range("I135").Value = 1022.0909
$ws.Range("K135").Value = 9198.8181
$ws.Range("M135").Value = -6663.8181
$ws.Range("H137").Value = 1841.7715
$ws.Range("I137").Value = 1396
$ws.Range("J137").Value = 2217.158
$ws.Range("K137").Value = 4188
$ws.Range("L137").Value = 6651.474
$ws.Range("M137").Value = -1638
$ws.Range("N137").Value = -11751.474
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5309.2373
$ws.Range("I32").Value = 2170.2134
$ws.Range("K32").Value = 2170.2134
$ws.Range("M32").Value = -1883.2134
$ws.Range("H45").Value = 1915.7894
$ws.Range("I45").Value = 1246.7142
$ws.Range("K45").Value = 1246.7142
$ws.Range("M45").Value = -869.7141999999999
$ws.Range("H61").Value = 8300.583000000001
$ws.Range("I61").Value = 8054.5454
$ws.Range("K61").Value = 8054.5454
$ws.Range("M61").Value = -7842.5454
$ws.Range("H74").Value = 4095.8
$ws.Range("I74").Value = 993.5333000000001
$ws.Range("J74").Value = 13402.6
$ws.Range("K74").Value = 993.5333000000001
$ws.Range("L74").Value = 13402.6
$ws.Range("M74").Value = -119.5333000000001
$ws.Range("N74").Value = -15150.6
$ws.Range("H77").Value = 4095.8
$ws.Range("I77").Value = 993.5333000000001
$ws.Range("J77").Value = 13402.6
$ws.Range("K77").Value = 4967.6665
$ws.Range("L77").Value = 67013
$ws.Range("M77").Value = -599.6665000000003
$ws.Range("N77").Value = -75749
$ws.Range("H82").Value = 151865
$ws.Range("J82").Value = 151865
$ws.Range("L82").Value = 151865
$ws.Range("N82").Value = -152587
$ws.Range("H85").Value = 151865
$ws.Range("J85").Value = 151865
$ws.Range("L85").Value = 151865
$ws.Range("N85").Value = -154361
$ws.Range("H122").Value = 3997.6667
$ws.Range("I122").Value = 3997.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11993.0001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9543.000100000001
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 8300.583000000001
$ws.Range("I136").Value = 8054.5454
$ws.Range("K136").Value = 24163.6362
$ws.Range("M136").Value = -21613.6362
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 4800
$ws.Range("I5").Value = 4950
$ws.Range("K5").Value = 4950
$ws.Range("M5").Value = -4837
$ws.Range("H20").Value = 1418.8823
$ws.Range("I20").Value = 1998.25
$ws.Range("J20").Value = 903.8889
$ws.Range("K20").Value = 1998.25
$ws.Range("L20").Value = 903.8889
$ws.Range("M20").Value = -1751.25
$ws.Range("N20").Value = -1397.8889
$ws.Range("H22").Value = 714.8570999999999
$ws.Range("I22").Value = 627.55554
$ws.Range("K22").Value = 627.55554
$ws.Range("M22").Value = -454.55554
$ws.Range("H54").Value = 15996.5
$ws.Range("I54").Value = 8040
$ws.Range("J54").Value = 18648.666
$ws.Range("K54").Value = 8040
$ws.Range("L54").Value = 18648.666
$ws.Range("M54").Value = -7556
$ws.Range("N54").Value = -19616.666
$ws.Range("H86").Value = 2183.0527
$ws.Range("I86").Value = 2193.4443
$ws.Range("J86").Value = 1996
$ws.Range("K86").Value = 2193.4443
$ws.Range("L86").Value = 1996
$ws.Range("M86").Value = -1070.4443
$ws.Range("N86").Value = -4242
$ws.Range("H89").Value = 2183.0527
$ws.Range("I89").Value = 2193.4443
$ws.Range("J89").Value = 1996
$ws.Range("K89").Value = 10967.2215
$ws.Range("L89").Value = 9980
$ws.Range("M89").Value = -5351.2215
$ws.Range("N89").Value = -21212
$ws.Range("H99").Value = 2455.4666
$ws.Range("I99").Value = 2179.4614
$ws.Range("K99").Value = 2179.4614
$ws.Range("M99").Value = -681.4614000000001
$ws.Range("H107").Value = 2078.4443
$ws.Range("I107").Value = 1402.7142
$ws.Range("K107").Value = 1402.7142
$ws.Range("M107").Value = 517.2858000000001
$ws.Range("H134").Value = 3206.7222
$ws.Range("I134").Value = 1951.3125
$ws.Range("K134").Value = 5853.9375
$ws.Range("M134").Value = -3318.9375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 687
$ws.Range("I22").Value = 582.8333
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 582.8333
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -232.8333
$ws.Range("N22").Value = -1699.5
$ws.Range("H33").Value = 2500
$ws.Range("I33").Value = 2500
$ws.Range("K33").Value = 2500
$ws.Range("M33").Value = -2121
$ws.Range("H62").Value = 4386.25
$ws.Range("I62").Value = 3205.125
$ws.Range("J62").Value = 6748.5
$ws.Range("K62").Value = 3205.125
$ws.Range("L62").Value = 6748.5
$ws.Range("M62").Value = -2581.125
$ws.Range("N62").Value = -7996.5
$ws.Range("H65").Value = 4386.25
$ws.Range("I65").Value = 3205.125
$ws.Range("J65").Value = 6748.5
$ws.Range("K65").Value = 16025.625
$ws.Range("L65").Value = 33742.5
$ws.Range("M65").Value = -12905.625
$ws.Range("N65").Value = -39982.5
$ws.Range("H80").Value = 153965.14
$ws.Range("J80").Value = 153965.14
$ws.Range("L80").Value = 153965.14
$ws.Range("N80").Value = -156211.14
$ws.Range("H83").Value = 153965.14
$ws.Range("J83").Value = 153965.14
$ws.Range("L83").Value = 461895.42
$ws.Range("N83").Value = -473127.42
$ws.Range("H94").Value = 2684
$ws.Range("I94").Value = 2087.7778
$ws.Range("K94").Value = 2087.7778
$ws.Range("M94").Value = -1636.7778
$ws.Range("H105").Value = 2419
$ws.Range("I105").Value = 1989.25
$ws.Range("K105").Value = 1989.25
$ws.Range("M105").Value = -242.25
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 1307
$ws.Range("I104").Value = 500
$ws.Range("K104").Value = 1500
$ws.Range("M104").Value = 1121
$ws.Range("H109").Value = 1135
$ws.Range("I109").Value = 1135
$ws.Range("K109").Value = 3405
$ws.Range("M109").Value = -2365
$ws.Range("H114").Value = 1715.875
$ws.Range("I114").Value = 1676
$ws.Range("J114").Value = 1739.8
$ws.Range("K114").Value = 5028
$ws.Range("L114").Value = 5219.4
$ws.Range("M114").Value = -1774
$ws.Range("N114").Value = -11727.4
$ws.Range("H128").Value = 449996.34
$ws.Range("I128").Value = 449996.34
$ws.Range("K128").Value = 1349989.02
$ws.Range("M128").Value = -1345009.02
$ws.Range("H129").Value = 8784188
$ws.Range("J129").Value = 15162405
$ws.Range("L129").Value = 45487215
$ws.Range("N129").Value = -45497215
$ws.Range("H131").Value = 2111.3914
$ws.Range("J131").Value = 4853.625
$ws.Range("L131").Value = 14560.875
$ws.Range("N131").Value = -24640.875
$ws.Range("H140").Value = 1782.4762
$ws.Range("I140").Value = 1766
$ws.Range("J140").Value = 1800.6
$ws.Range("K140").Value = 5298
$ws.Range("L140").Value = 5401.799999999999
$ws.Range("M140").Value = -118
$ws.Range("N140").Value = -15761.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 199500
$ws.Range("J15").Value = 199500
$ws.Range("L15").Value = 199500
$ws.Range("N15").Value = -200076
$ws.Range("H69").Value = 199500
$ws.Range("J69").Value = 199500
$ws.Range("L69").Value = 199500
$ws.Range("N69").Value = -200998
$ws.Range("H70").Value = 5540
$ws.Range("I70").Value = 5112.5
$ws.Range("K70").Value = 5112.5
$ws.Range("M70").Value = -4842.5
$ws.Range("H72").Value = 199500
$ws.Range("J72").Value = 199500
$ws.Range("L72").Value = 598500
$ws.Range("N72").Value = -605988
$ws.Range("H73").Value = 5540
$ws.Range("I73").Value = 5112.5
$ws.Range("K73").Value = 5112.5
$ws.Range("M73").Value = -4176.5
$ws.Range("H74").Value = 199500
$ws.Range("J74").Value = 199500
$ws.Range("L74").Value = 199500
$ws.Range("N74").Value = -201372
$ws.Range("H77").Value = 199500
$ws.Range("J77").Value = 199500
$ws.Range("L77").Value = 598500
$ws.Range("N77").Value = -607860
$ws.Range("H80").Value = 12074.125
$ws.Range("I80").Value = 10461.125
$ws.Range("K80").Value = 10461.125
$ws.Range("M80").Value = -9463.125
$ws.Range("H81").Value = 199500
$ws.Range("J81").Value = 199500
$ws.Range("L81").Value = 199500
$ws.Range("N81").Value = -201496
$ws.Range("H83").Value = 12074.125
$ws.Range("I83").Value = 10461.125
$ws.Range("K83").Value = 52305.625
$ws.Range("M83").Value = -47313.625
$ws.Range("H84").Value = 199500
$ws.Range("J84").Value = 199500
$ws.Range("L84").Value = 598500
$ws.Range("N84").Value = -608484
$ws.Range("H87").Value = 199500
$ws.Range("J87").Value = 199500
$ws.Range("L87").Value = 199500
$ws.Range("N87").Value = -201996
$ws.Range("H90").Value = 199500
$ws.Range("J90").Value = 199500
$ws.Range("L90").Value = 598500
$ws.Range("N90").Value = -610980
$ws.Range("H97").Value = 452.9524
$ws.Range("J97").Value = 740.8570999999999
$ws.Range("L97").Value = 740.8570999999999
$ws.Range("N97").Value = -1732.8571
$ws.Range("H107").Value = 657.7727
$ws.Range("I107").Value = 482.84616
$ws.Range("K107").Value = 482.84616
$ws.Range("M107").Value = 1437.15384
$ws.Range("H122").Value = 3822.52
$ws.Range("I122").Value = 1418.5238
$ws.Range("J122").Value = 16443.5
$ws.Range("K122").Value = 4255.5714
$ws.Range("L122").Value = 49330.5
$ws.Range("M122").Value = -1805.5714
$ws.Range("N122").Value = -54230.5
$ws.Range("H126").Value = 3776.36
$ws.Range("J126").Value = 5788.364
$ws.Range("L126").Value = 17365.092
$ws.Range("N126").Value = -22305.092
$ws.Range("H132").Value = 4898.364
$ws.Range("I132").Value = 2927.8333
$ws.Range("J132").Value = 7263
$ws.Range("K132").Value = 8783.499899999999
$ws.Range("L132").Value = 21789
$ws.Range("M132").Value = -6253.499899999999
$ws.Range("N132").Value = -26849
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 54580.5
$ws.Range("I7").Value = 58088.65
$ws.Range("K7").Value = 58088.65
$ws.Range("M7").Value = -57976.65
$ws.Range("H22").Value = 2649.25
$ws.Range("H26").Value = 30000
$ws.Range("I26").Value = 30000
$ws.Range("K26").Value = 30000
$ws.Range("M26").Value = -29705
$ws.Range("H27").Value = 2649.25
$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 20000
$ws.Range("K29").Value = 20000
$ws.Range("M29").Value = -19705
$ws.Range("H40").Value = 4909.5
$ws.Range("I40").Value = 3883.7856
$ws.Range("K40").Value = 3883.7856
$ws.Range("M40").Value = -3747.7856
$ws.Range("H46").Value = 2734.1365
$ws.Range("I46").Value = 640.5
$ws.Range("J46").Value = 3519.25
$ws.Range("K46").Value = 640.5
$ws.Range("L46").Value = 3519.25
$ws.Range("M46").Value = -452.5
$ws.Range("N46").Value = -3895.25
$ws.Range("H69").Value = 199500
$ws.Range("J69").Value = 199500
$ws.Range("L69").Value = 199500
$ws.Range("N69").Value = -201122
$ws.Range("H72").Value = 199500
$ws.Range("J72").Value = 199500
$ws.Range("L72").Value = 598500
$ws.Range("N72").Value = -606612
$ws.Range("H80").Value = 95029.414
$ws.Range("J80").Value = 95029.414
$ws.Range("L80").Value = 95029.414
$ws.Range("N80").Value = -97275.414
$ws.Range("H82").Value = 1372.3
$ws.Range("I82").Value = 1065
$ws.Range("J82").Value = 1504
$ws.Range("K82").Value = 1065
$ws.Range("L82").Value = 1504
$ws.Range("M82").Value = -704
$ws.Range("N82").Value = -2226
$ws.Range("H83").Value = 95029.414
$ws.Range("J83").Value = 95029.414
$ws.Range("L83").Value = 285088.242
$ws.Range("N83").Value = -296320.242
$ws.Range("H85").Value = 1372.3
$ws.Range("I85").Value = 1065
$ws.Range("J85").Value = 1504
$ws.Range("K85").Value = 1065
$ws.Range("L85").Value = 1504
$ws.Range("M85").Value = 183
$ws.Range("N85").Value = -4000
$ws.Range("H100").Value = 1450.1818
$ws.Range("I100").Value = 1345.2
$ws.Range("K100").Value = 1345.2
$ws.Range("M100").Value = -804.2
$ws.Range("H122").Value = 3908.3438
$ws.Range("I122").Value = 3223.6667
$ws.Range("K122").Value = 9671.000100000001
$ws.Range("M122").Value = -7221.000100000001
$ws.Range("H126").Value = 54580.5
$ws.Range("I126").Value = 58088.65
$ws.Range("K126").Value = 174265.95
$ws.Range("M126").Value = -171795.95
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 12502
$ws.Range("I132").Value = 7004
$ws.Range("K132").Value = 21012
$ws.Range("M132").Value = -18482
$ws.Range("H136").Value = 7074.625
$ws.Range("I136").Value = 4985.2856
$ws.Range("J136").Value = 8699.666999999999
$ws.Range("K136").Value = 14955.8568
$ws.Range("L136").Value = 26099.001
$ws.Range("M136").Value = -12405.8568
$ws.Range("N136").Value = -31199.001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 4949
$ws.Range("J8").Value = 4949
$ws.Range("L8").Value = 4949
$ws.Range("N8").Value = -5229
$ws.Range("H62").Value = 5536.1665
$ws.Range("I62").Value = 5928.6665
$ws.Range("J62").Value = 5143.6665
$ws.Range("K62").Value = 5928.6665
$ws.Range("L62").Value = 5143.6665
$ws.Range("M62").Value = -5304.6665
$ws.Range("N62").Value = -6391.6665
$ws.Range("H65").Value = 5536.1665
$ws.Range("I65").Value = 5928.6665
$ws.Range("J65").Value = 5143.6665
$ws.Range("K65").Value = 29643.3325
$ws.Range("L65").Value = 25718.3325
$ws.Range("M65").Value = -26523.3325
$ws.Range("N65").Value = -31958.3325
$ws.Range("H76").Value = 152117.58
$ws.Range("J76").Value = 152117.58
$ws.Range("L76").Value = 152117.58
$ws.Range("N76").Value = -152747.58
$ws.Range("H79").Value = 152117.58
$ws.Range("J79").Value = 152117.58
$ws.Range("L79").Value = 152117.58
$ws.Range("N79").Value = -154301.58
$ws.Range("H81").Value = 2793
$ws.Range("I81").Value = 2793
$ws.Range("K81").Value = 5586
$ws.Range("M81").Value = -4525
$ws.Range("H84").Value = 2793
$ws.Range("I84").Value = 2793
$ws.Range("K84").Value = 27930
$ws.Range("M84").Value = -22626
$ws.Range("H87").Value = 199500
$ws.Range("J87").Value = 199500
$ws.Range("L87").Value = 199500
$ws.Range("N87").Value = -201996
$ws.Range("H90").Value = 199500
$ws.Range("J90").Value = 199500
$ws.Range("L90").Value = 598500
$ws.Range("N90").Value = -610980
$ws.Range("H100").Value = 389.25
$ws.Range("I100").Value = 389.25
$ws.Range("K100").Value = 778.5
$ws.Range("M100").Value = -237.5
$ws.Range("H107").Value = 741452.5
$ws.Range("I107").Value = 1000606.9
$ws.Range("J107").Value = 1011.4286
$ws.Range("K107").Value = 3001820.7
$ws.Range("L107").Value = 3034.2858
$ws.Range("M107").Value = -2999900.7
$ws.Range("N107").Value = -6874.2858
$ws.Range("H122").Value = 2751.913
$ws.Range("J122").Value = 4000
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 1982.4517
$ws.Range("I132").Value = 1423.3448
$ws.Range("J132").Value = 10089.5
$ws.Range("K132").Value = 4270.0344
$ws.Range("L132").Value = 30268.5
$ws.Range("M132").Value = -1740.0344
$ws.Range("N132").Value = -35328.5
$ws.Range("H136").Value = 5726.8613
$ws.Range("I136").Value = 5488.25
$ws.Range("J136").Value = 7635.75
$ws.Range("K136").Value = 16464.75
$ws.Range("L136").Value = 22907.25
$ws.Range("M136").Value = -13914.75
$ws.Range("N136").Value = -28007.25
